$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the col_value cells as quote-prefixed text so the leading zeros
# (and the literal surrounding apostrophes) are preserved, matching how
# the source Python script re-wrote these values.
$ws.Range("C3").Value = "''253'"
$ws.Range("C4").Value = "''0010'"
$ws.Range("C5").Value = "''0000'"

# Move the active selection to C6 (matches the saved cursor position).
$ws.Range("C6").Select()
